# Update computed market-price / profit columns (H:N) on several Leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the latest
# scheduled market-data pull. Only numeric columns H..N change; the
# Leve metadata columns (A..G) are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Cells.Item(4, 8).Value = 183.16667
$ws.Cells.Item(4, 9).Value = 149.8
$ws.Cells.Item(4, 10).Value = 350
$ws.Cells.Item(4, 11).Value = 149.8
$ws.Cells.Item(4, 12).Value = 350
$ws.Cells.Item(4, 13).Value = -35.80000000000001
$ws.Cells.Item(4, 14).Value = -578
# Row 9
$ws.Cells.Item(9, 8).Value = 959.8461
$ws.Cells.Item(9, 9).Value = 1039.091
$ws.Cells.Item(9, 10).Value = 524
$ws.Cells.Item(9, 11).Value = 1039.091
$ws.Cells.Item(9, 12).Value = 524
$ws.Cells.Item(9, 13).Value = -870.0909999999999
$ws.Cells.Item(9, 14).Value = -862
# Row 13
$ws.Cells.Item(13, 8).Value = 6499.5
$ws.Cells.Item(13, 9).Value = 8000
$ws.Cells.Item(13, 10).Value = 4999
$ws.Cells.Item(13, 11).Value = 8000
$ws.Cells.Item(13, 12).Value = 4999
$ws.Cells.Item(13, 13).Value = -7831
$ws.Cells.Item(13, 14).Value = -5337
# Row 33
$ws.Cells.Item(33, 8).Value = 331.4
$ws.Cells.Item(33, 9).Value = 334.94446
$ws.Cells.Item(33, 11).Value = 334.94446
$ws.Cells.Item(33, 13).Value = -105.94446
# Row 41
$ws.Cells.Item(41, 8).Value = 290.57144
$ws.Cells.Item(41, 9).Value = 302
$ws.Cells.Item(41, 11).Value = 302
$ws.Cells.Item(41, 13).Value = 138
# Row 62
$ws.Cells.Item(62, 8).Value = 4495.5
$ws.Cells.Item(62, 9).Value = 4489
$ws.Cells.Item(62, 11).Value = 4489
$ws.Cells.Item(62, 13).Value = -3865
# Row 65
$ws.Cells.Item(65, 8).Value = 4495.5
$ws.Cells.Item(65, 9).Value = 4489
$ws.Cells.Item(65, 11).Value = 22445
$ws.Cells.Item(65, 13).Value = -19325
# Row 92
$ws.Cells.Item(92, 8).Value = 40719.48
$ws.Cells.Item(92, 9).Value = 59472.883
$ws.Cells.Item(92, 10).Value = 868.5
$ws.Cells.Item(92, 11).Value = 59472.883
$ws.Cells.Item(92, 12).Value = 868.5
$ws.Cells.Item(92, 13).Value = -58224.883
$ws.Cells.Item(92, 14).Value = -3364.5
# Row 137
$ws.Cells.Item(137, 8).Value = 4466.5
$ws.Cells.Item(137, 9).Value = 1800
$ws.Cells.Item(137, 10).Value = 4999.8
$ws.Cells.Item(137, 11).Value = 5400
$ws.Cells.Item(137, 12).Value = 14999.4
$ws.Cells.Item(137, 13).Value = -2850
$ws.Cells.Item(137, 14).Value = -20099.4

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 10943.929
$ws.Cells.Item(32, 9).Value = 10460.223
$ws.Cells.Item(32, 11).Value = 10460.223
$ws.Cells.Item(32, 13).Value = -10173.223
# Row 74
$ws.Cells.Item(74, 8).Value = 1788.5883
$ws.Cells.Item(74, 9).Value = 1818.5
$ws.Cells.Item(74, 11).Value = 1818.5
$ws.Cells.Item(74, 13).Value = -944.5
# Row 77
$ws.Cells.Item(77, 8).Value = 1788.5883
$ws.Cells.Item(77, 9).Value = 1818.5
$ws.Cells.Item(77, 11).Value = 9092.5
$ws.Cells.Item(77, 13).Value = -4724.5
# Row 80
$ws.Cells.Item(80, 8).Value = 54399.6
$ws.Cells.Item(80, 10).Value = 149999
$ws.Cells.Item(80, 12).Value = 149999
$ws.Cells.Item(80, 14).Value = -151995
# Row 83
$ws.Cells.Item(83, 8).Value = 54399.6
$ws.Cells.Item(83, 10).Value = 149999
$ws.Cells.Item(83, 12).Value = 449997
$ws.Cells.Item(83, 14).Value = -459981

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Cells.Item(134, 8).Value = 4615.0386
$ws.Cells.Item(134, 9).Value = 4443.24
$ws.Cells.Item(134, 11).Value = 13329.72
$ws.Cells.Item(134, 13).Value = -10794.72

$ws = $wb.Worksheets.Item("CRP")
# Row 55
$ws.Cells.Item(55, 8).Value = 0
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 12).Value = 0
$ws.Cells.Item(55, 14).ClearContents()  # N55: was -13630, now blank
# Row 58
$ws.Cells.Item(58, 8).Value = 2655.5557
$ws.Cells.Item(58, 9).Value = 1435.9286
$ws.Cells.Item(58, 11).Value = 1435.9286
$ws.Cells.Item(58, 13).Value = -1232.9286
# Row 132
$ws.Cells.Item(132, 8).Value = 1922.1428
$ws.Cells.Item(132, 9).Value = 1839.5385
$ws.Cells.Item(132, 10).Value = 2996
$ws.Cells.Item(132, 11).Value = 5518.6155
$ws.Cells.Item(132, 12).Value = 8988
$ws.Cells.Item(132, 13).Value = -2988.6155
$ws.Cells.Item(132, 14).Value = -14048
# Row 136
$ws.Cells.Item(136, 8).Value = 2655.5557
$ws.Cells.Item(136, 9).Value = 1435.9286
$ws.Cells.Item(136, 11).Value = 4307.7858
$ws.Cells.Item(136, 13).Value = -1757.7858

$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Cells.Item(55, 8).Value = 2978.375
$ws.Cells.Item(55, 10).Value = 3645.1538
$ws.Cells.Item(55, 12).Value = 10935.4614
$ws.Cells.Item(55, 14).Value = -11289.4614
# Row 68
$ws.Cells.Item(68, 8).Value = 999
$ws.Cells.Item(68, 9).Value = 999
$ws.Cells.Item(68, 10).Value = 999
$ws.Cells.Item(68, 11).Value = 2997
$ws.Cells.Item(68, 12).Value = 2997
$ws.Cells.Item(68, 13).Value = -2186
$ws.Cells.Item(68, 14).Value = -4619
# Row 71
$ws.Cells.Item(71, 8).Value = 999
$ws.Cells.Item(71, 9).Value = 999
$ws.Cells.Item(71, 10).Value = 999
$ws.Cells.Item(71, 11).Value = 8991
$ws.Cells.Item(71, 12).Value = 8991
$ws.Cells.Item(71, 13).Value = -4935
$ws.Cells.Item(71, 14).Value = -17103
# Row 97
$ws.Cells.Item(97, 8).Value = 1528
$ws.Cells.Item(97, 9).Value = 1534.8
$ws.Cells.Item(97, 10).Value = 1494
$ws.Cells.Item(97, 11).Value = 4604.4
$ws.Cells.Item(97, 12).Value = 4482
$ws.Cells.Item(97, 13).Value = -4108.4
$ws.Cells.Item(97, 14).Value = -5474
# Row 132
$ws.Cells.Item(132, 8).Value = 1748.1666
$ws.Cells.Item(132, 9).Value = 1747.5
$ws.Cells.Item(132, 10).Value = 1748.5
$ws.Cells.Item(132, 11).Value = 15727.5
$ws.Cells.Item(132, 12).Value = 15736.5
$ws.Cells.Item(132, 13).Value = -13197.5
$ws.Cells.Item(132, 14).Value = -20796.5

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 9042.076999999999
$ws.Cells.Item(70, 9).Value = 8817.75
$ws.Cells.Item(70, 11).Value = 8817.75
$ws.Cells.Item(70, 13).Value = -8547.75
# Row 73
$ws.Cells.Item(73, 8).Value = 9042.076999999999
$ws.Cells.Item(73, 9).Value = 8817.75
$ws.Cells.Item(73, 11).Value = 8817.75
$ws.Cells.Item(73, 13).Value = -7881.75

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Cells.Item(46, 8).Value = 1173
$ws.Cells.Item(46, 9).Value = 1222.625
$ws.Cells.Item(46, 10).Value = 974.5
$ws.Cells.Item(46, 11).Value = 1222.625
$ws.Cells.Item(46, 12).Value = 974.5
$ws.Cells.Item(46, 13).Value = -1034.625
$ws.Cells.Item(46, 14).Value = -1350.5
# Row 93
$ws.Cells.Item(93, 8).Value = 1125.3334
$ws.Cells.Item(93, 9).Value = 1328.0667
$ws.Cells.Item(93, 10).Value = 618.5
$ws.Cells.Item(93, 11).Value = 1328.0667
$ws.Cells.Item(93, 12).Value = 618.5
$ws.Cells.Item(93, 13).Value = -80.06670000000008
$ws.Cells.Item(93, 14).Value = -3114.5
# Row 100
$ws.Cells.Item(100, 8).Value = 2010.375
$ws.Cells.Item(100, 9).Value = 1997.5714
$ws.Cells.Item(100, 11).Value = 1997.5714
$ws.Cells.Item(100, 13).Value = -1456.5714
# Row 132
$ws.Cells.Item(132, 8).Value = 2156.606
$ws.Cells.Item(132, 9).Value = 2346.88
$ws.Cells.Item(132, 11).Value = 7040.64
$ws.Cells.Item(132, 13).Value = -4510.64
# Row 136
$ws.Cells.Item(136, 8).Value = 5806.16
$ws.Cells.Item(136, 9).Value = 4507.85
$ws.Cells.Item(136, 10).Value = 10999.4
$ws.Cells.Item(136, 11).Value = 13523.55
$ws.Cells.Item(136, 12).Value = 32998.2
$ws.Cells.Item(136, 13).Value = -10973.55
$ws.Cells.Item(136, 14).Value = -38098.2
# Row 140
$ws.Cells.Item(140, 8).Value = 88429
$ws.Cells.Item(140, 10).Value = 88429
$ws.Cells.Item(140, 12).Value = 88429
$ws.Cells.Item(140, 14).Value = -98789

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Cells.Item(5, 8).Value = 1500
$ws.Cells.Item(5, 10).Value = 1500
$ws.Cells.Item(5, 12).Value = 1500
$ws.Cells.Item(5, 14).Value = -1724
# Row 113
$ws.Cells.Item(113, 8).Value = 601.9286
$ws.Cells.Item(113, 9).Value = 453.8
$ws.Cells.Item(113, 10).Value = 972.25
$ws.Cells.Item(113, 11).Value = 1361.4
$ws.Cells.Item(113, 12).Value = 2916.75
$ws.Cells.Item(113, 13).Value = 808.5999999999999
$ws.Cells.Item(113, 14).Value = -7256.75
# Row 122
$ws.Cells.Item(122, 8).Value = 6677.923
$ws.Cells.Item(122, 9).Value = 5900.6665
$ws.Cells.Item(122, 10).Value = 16005
$ws.Cells.Item(122, 11).Value = 17701.9995
$ws.Cells.Item(122, 12).Value = 48015
$ws.Cells.Item(122, 13).Value = -15251.9995
$ws.Cells.Item(122, 14).Value = -52915
